$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new timesheet entry row 9 (copy A8's formatting so the date style is reused)
$ws.Range("A8").Copy($ws.Range("A9"))
$ws.Range("A9").Value = 43381
$ws.Range("B9").Value = 0.75
$ws.Range("C9").Value = "Coderen fietssimulatie"

# Update the selection to match the new active cell
$ws.Range("C10").Select()
